$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.671.25'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '2.500.98'
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'575.86"
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = "'166.98"
$ws.Range("E6").Value = '  +0.54%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").Value = '2.500.32'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  +1.65%  '
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = '  +3.43%  '
$ws.Range("D13").Value = "'4.96"
$ws.Range("E13").Value = '  +2.19%  '
$ws.Range("D14").Value = '2.952.21'
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("D15").Value = '69.617.84'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = '  +2.53%  '
$ws.Range("D17").Value = "'24.77"
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("D18").Value = '2.495.49'
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").Value = "'11.19"
$ws.Range("E19").Value = '  -1.60%  '
$ws.Range("D20").Value = "'7.46"
$ws.Range("E20").Value = '  -4.47%  '
$ws.Range("D21").Value = "'348.55"
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").Value = "'3.90"
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = "'70.39"
$ws.Range("E25").Value = '  +2.70%  '
$ws.Range("D26").Value = "'3.96"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = "'8.77"
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("D28").Value = '2.605.88'
$ws.Range("E28").Value = '  -0.94%  '
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '0.0₃0891'
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("D31").Value = "'7.82"
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").Value = "'459.91"
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("E33").Value = '  -2.76%  '
$ws.Range("D34").Value = "'1.73"
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = "'156.67"
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("D39").Value = "'18.49"
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("D42").Value = "'4.68"
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = "'38.12"
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = "'2.20"
$ws.Range("E45").Value = '  -4.00%  '
$ws.Range("E46").Value = '  -6.46%  '
$ws.Range("D47").Value = "'141.16"
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = "'0.518"
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("D50").Value = "'0.0733"
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").Value = "'0.578"
$ws.Range("E51").Value = '  -0.95%  '
